# Commit: "Mon, Jun 15, 2020  1:05:46 PM"
#
# 1) Slide 16 table: change the table style to the built-in
#    {12D3D91C-9AA8-4FBE-B400-81104440F2BA} style.
# 2) Theme: recolor the presentation's main theme (theme1.xml, used by
#    the slide master) from the "Integral" palette to the standard
#    "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 16 -------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{12D3D91C-9AA8-4FBE-B400-81104440F2BA}")
    }
}

# --- 2) Theme colours -------------------------------------------------------
$master = $p.SlideMaster
$colors = $master.ColorScheme

# Office Theme palette, in ColorScheme.Colors(1..12) order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
# accent6, hlink, folHlink  (RGB encoded as R + G*256 + B*65536)
$colors.Colors(1).RGB  = 0          # dk1      000000
$colors.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Colors(3).RGB  = 6968388    # dk2      44546A
$colors.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Colors(6).RGB  = 3243501    # accent2  ED7D31
$colors.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Colors(8).RGB  = 49407      # accent4  FFC000
$colors.Colors(9).RGB  = 12874308   # accent5  4472C4
$colors.Colors(10).RGB = 4697456    # accent6  70AD47
$colors.Colors(11).RGB = 12673797   # hlink    0563C1
$colors.Colors(12).RGB = 7491477    # folHlink 954F72
